$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L5").Value = 940.62
$ws1.Range("M5").Value = 12302.59
$ws1.Range("D11").Value = 648

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 13243.21
$ws2.Range("F11").Value = 6856.03
$ws2.Range("F23").Value = 56313.78

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths adjusted (D: 14 -> 13, E: 23 -> 24)
# Note: the runtime's ColumnWidth setter applies a fixed +5/6 (0.8333..) offset
# versus the stored OOXML column width, so compensate by subtracting it here
# to land exactly on the target stored widths of 13 and 24.
$ws3.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668

$ws3.Range("D3").Value = 2145.63
$ws3.Range("E3").Value = 3358.98890386263
$ws3.Range("F3").Value = 0.3897872018886532

$ws3.Range("D11").Value = 3439.54
$ws3.Range("E11").Value = 2404.90916370549
$ws3.Range("F11").Value = 0.5885139734570414

$ws3.Range("D12").Value = 47392.3
$ws3.Range("E12").Value = -9652.560000000005
$ws3.Range("F12").Value = 1.255766467919493

$ws3.Range("D14").Value = 56313.78
$ws3.Range("E14").Value = -889.0385211961108
$ws3.Range("F14").Value = 1.016040463112239
